# Update "想去人数" (want-to-go count) values in column F across the
# four worksheets of the 广州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (sheet1) ---
$ws1.Range("F3").Value = 650
$ws1.Range("F4").Value = 811
$ws1.Range("F5").Value = 497
$ws1.Range("F6").Value = 389
$ws1.Range("F7").Value = 476
$ws1.Range("F10").Value = 816
$ws1.Range("F11").Value = 640
$ws1.Range("F12").Value = 113
$ws1.Range("F15").Value = 717
$ws1.Range("F16").Value = 210
$ws1.Range("F17").Value = 513
$ws1.Range("F18").Value = 461
$ws1.Range("F19").Value = 1227
$ws1.Range("F21").Value = 928
$ws1.Range("F22").Value = 2678
$ws1.Range("F23").Value = 1150
$ws1.Range("F24").Value = 611
$ws1.Range("F25").Value = 138
$ws1.Range("F26").Value = 1196
$ws1.Range("F27").Value = 49
$ws1.Range("F28").Value = 890
$ws1.Range("F29").Value = 90
$ws1.Range("F30").Value = 1224

# --- Sheet "演出" (sheet2) ---
$ws2.Range("F3").Value = 478

# --- Sheet "本地生活" (sheet3) ---
$ws3.Range("F2").Value = 696

# --- Sheet "全部类型" (sheet4) ---
$ws4.Range("F2").Value = 696
$ws4.Range("F4").Value = 650
$ws4.Range("F5").Value = 811
$ws4.Range("F6").Value = 497
$ws4.Range("F8").Value = 389
$ws4.Range("F9").Value = 476
$ws4.Range("F10").Value = 478
$ws4.Range("F16").Value = 816
$ws4.Range("F17").Value = 640
$ws4.Range("F18").Value = 113
$ws4.Range("F26").Value = 717
$ws4.Range("F27").Value = 210
$ws4.Range("F28").Value = 513
$ws4.Range("F29").Value = 461
$ws4.Range("F30").Value = 1227
$ws4.Range("F32").Value = 928
$ws4.Range("F33").Value = 2678
$ws4.Range("F34").Value = 1150
$ws4.Range("F35").Value = 611
$ws4.Range("F36").Value = 138
$ws4.Range("F37").Value = 1196
$ws4.Range("F38").Value = 49
$ws4.Range("F40").Value = 890
$ws4.Range("F41").Value = 90
$ws4.Range("F42").Value = 1224
